$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NonDurable")
$ws.Range("B1:C16").EntireColumn.AutoFit()
